# Applies the "moving files, cleaning indiv no file" edit:
# fills in a number of previously-blank "y" observation cells (shared
# string index 25 == "y") across the GC_splist sheet, and updates the
# active selection/zoom to match the author's final view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells that already carry the correct fill/style (s="1") and only
#     need their value set to "y". Setting .Value preserves the cell's
#     existing style index, matching the diff exactly.
$valueOnlyCells = @("E4", "F4", "G4", "G8", "G9", "G10", "G16")
foreach ($addr in $valueOnlyCells) {
    $ws.Range($addr).Value = "y"
}

# --- Brand new cells that need no special fill (style 0) -- just set
#     the value; Excel will create a plain cell.
$plainNewCells = @("H2", "H3", "H7", "H11", "I11", "G12", "I13", "H14", "H15", "H16", "H18", "H21", "G22")
foreach ($addr in $plainNewCells) {
    $ws.Range($addr).Value = "y"
}

# --- Brand new cells that need to pick up the same shaded fill style
#     used by neighbouring cells (style index 1). Copy the format from
#     a same-row cell that already carries that style, then set value.
$styledNewCells = @(
    @{Target="H4"; Source="D4"},
    @{Target="I4"; Source="D4"},
    @{Target="H6"; Source="G6"},
    @{Target="I6"; Source="G6"},
    @{Target="H8"; Source="F8"}
)
foreach ($pair in $styledNewCells) {
    $ws.Range($pair.Source).Copy()
    $ws.Range($pair.Target).PasteSpecial(-4122)
    $ws.Range($pair.Target).Value = "y"
}

$excel.CutCopyMode = 0

# --- Update the saved view: zoom level and active selection.
$excel.ActiveWindow.Zoom = 140
$ws.Range("H15").Select()
